$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Campaign Funds" column (D)
$ws.Range("D1").Value = "Campaign Funds"
$ws.Range("D5").Value = 4552
$ws.Range("D6").Value = 8442.01
$ws.Range("D7").Value = 6761
$ws.Range("D8").Value = 2800
$ws.Range("D9").Value = 20721
$ws.Range("D10").Value = 305
$ws.Range("D11").Value = 3127.04
$ws.Range("D12").Value = 5
$ws.Range("D13").Value = 5163

# Move the active selection to D14
[void]$ws.Range("D14").Select()
